# Apply updated crypto price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number (e.g. "207.10")
# must be forced to remain plain text, matching the original inlineStr cells,
# otherwise Excel auto-converts them to floating point numbers.
$ws.Range("D2").Value = "27.657.37"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.587.70"
$ws.Range("E3").Value = "  -2.61%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.10"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.499"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.58%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.21"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.253"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("E10").Value = "  -2.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0866"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").Value = "1.813.09"
$ws.Range("E12").Value = "  -2.57%  "
$ws.Range("D13").Value = "1.594.01"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("E14").Value = "  -4.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.530"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.81%  "
$ws.Range("D16").Value = "27.646.66"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.42"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "219.17"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.27%  "
$ws.Range("D19").Value = "0.0₃0695"
$ws.Range("E19").Value = "  -3.28%  "
$ws.Range("E20").Value = "  -4.24%  "
$ws.Range("E22").Value = "  -4.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.64"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.98"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.51"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.82"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.11"
$ws.Range("D28").ClearFormats()
$ws.Range("E29").Value = "  -4.92%  "
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("E32").Value = "  -5.55%  "
$ws.Range("D33").Value = "1.367.36"
$ws.Range("E33").Value = "  -3.57%  "
$ws.Range("E34").Value = "  -5.66%  "
$ws.Range("E35").Value = "  -4.95%  "
$ws.Range("E36").Value = "  -3.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.31"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0168"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.534"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.824"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.57%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("E42").Value = "  -3.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.11"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.16"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.12%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.724.20"
$ws.Range("E45").Value = "  -2.58%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.16"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.78%  "
$ws.Range("E47").Value = "  -5.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.57"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("D49").Value = "0.0₆0101"
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0964"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.83%  "
$ws.Range("E51").Value = "  -1.51%  "
